$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '42.706.79'
Set-TextValue "D3" '2.309.45'
Set-TextValue "E3" '  +0.61%  '
Set-TextValue "E4" '  +0.01%  '
Set-TextValue "D5" '301.58'
Set-TextValue "E5" '  -1.13%  '
Set-TextValue "D6" '95.33'
Set-TextValue "E6" '  -1.09%  '
Set-TextValue "D7" '0.502'
Set-TextValue "E7" '  -0.43%  '
Set-TextValue "E8" '  +0.02%  '
Set-TextValue "E9" '  -1.26%  '
Set-TextValue "E10" '  -2.67%  '
Set-TextValue "D11" '18.84'
Set-TextValue "E11" '  +1.30%  '
Set-TextValue "E12" '  -0.05%  '
Set-TextValue "E13" '  -0.04%  '
Set-TextValue "E14" '  -1.84%  '
Set-TextValue "D15" '2.670.91'
Set-TextValue "E15" '  +0.68%  '
Set-TextValue "D16" '2.279.95'
Set-TextValue "E16" '  +0.18%  '
Set-TextValue "D17" '0.788'
Set-TextValue "E17" '  +1.42%  '
Set-TextValue "D18" '42.656.18'
Set-TextValue "D19" '12.13'
Set-TextValue "E19" '  -4.76%  '
Set-TextValue "E20" '  +1.72%  '
Set-TextValue "D21" '0.0₃0889'
Set-TextValue "E22" '  +0.84%  '
Set-TextValue "D23" '2.26'
Set-TextValue "E23" '  +4.94%  '
Set-TextValue "D24" '235.33'
Set-TextValue "E24" '  -0.11%  '
Set-TextValue "E25" '  +0.02%  '
Set-TextValue "D26" '2.42'
Set-TextValue "E26" '  +0.55%  '
Set-TextValue "D27" '24.25'
Set-TextValue "E27" '  -1.60%  '
Set-TextValue "D28" '2.37'
Set-TextValue "E28" '  +15.54%  '
Set-TextValue "D29" '166.48'
Set-TextValue "D30" '9.09'
Set-TextValue "E30" '  +0.83%  '
Set-TextValue "D31" '32.08'
Set-TextValue "E31" '  -2.17%  '
Set-TextValue "E32" '  -0.03%  '
Set-TextValue "E33" '  +0.57%  '
Set-TextValue "D34" '17.62'
Set-TextValue "E34" '  -2.28%  '
Set-TextValue "D35" '4.46'
Set-TextValue "E35" '  +0.87%  '
Set-TextValue "D36" '0.0697'
Set-TextValue "E37" '  -0.97%  '
Set-TextValue "E38" '  +2.33%  '
Set-TextValue "E39" '  -0.12%  '
Set-TextValue "E40" '  -0.75%  '
Set-TextValue "D41" '2.69'
Set-TextValue "E41" '  -0.07%  '
Set-TextValue "D42" '21.15'
Set-TextValue "E42" '  +16.45%  '
Set-TextValue "D43" '1.922.88'
Set-TextValue "E43" '  -3.69%  '
Set-TextValue "E44" '  -0.48%  '
Set-TextValue "D45" '10.03'
Set-TextValue "E45" '  -2.20%  '
Set-TextValue "E46" '  -1.58%  '
Set-TextValue "E47" '  -1.07%  '
Set-TextValue "E48" '  +1.51%  '
Set-TextValue "D49" '2.539.53'
Set-TextValue "E49" '  +0.77%  '
Set-TextValue "D50" '53.28'
Set-TextValue "E50" '  -0.15%  '
Set-TextValue "D51" '72.02'
Set-TextValue "E51" '  +1.65%  '
